# Adds a new "IMAGENESU" table (idimagenesU / imagen) to the little table
# catalogue on Hoja1, re-orders TIPOS/MATERIAS/AULAS so IMAGENESU and TIPOS
# sit right after USUARIOS, and extends the USUARIOS column with the new
# idtipo / idimagen / nombre / apellido / contraseña / mail fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new IMAGENESU column: insert a blank column at B.
# This shifts MATERIAS/AULAS/TIPOS (old B:D) to C:E and MATERIAS-AULAS /
# MATERIAS-USUARIOS / TURNOS (old E:G, which carry explicit bestFit column
# widths) to F:H, preserving their widths automatically.
$ws.Columns.Item(2).Insert()

# --- New IMAGENESU table header/id (col B) and reordered TIPOS (col C) ---
$ws.Range("A5").Value = "idimagen"
$ws.Range("B2").Value = "IMAGENESU"
$ws.Range("C2").Value = "TIPOS"
$ws.Range("D2").Value = "MATERIAS"
$ws.Range("E2").Value = "AULAS"

$ws.Range("B3").Value = "idimagenesU"
$ws.Range("C3").Value = "idtipo"
$ws.Range("D3").Value = "idmateria"
$ws.Range("E3").Value = "idaula"

$ws.Range("A4").Value = "idtipo"
$ws.Range("B4").Value = "imagen"

$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = "idturno"

# --- Extra USUARIOS fields appended below (firestore/jwt columns) ---
$ws.Range("A6").Value = "nombre"
$ws.Range("A7").Value = "apellido"
$ws.Range("A8").Value = "contraseña"
$ws.Range("A9").Value = "mail"

# Restore the original selection location noted in the saved file.
$ws.Range("C2").Select() | Out-Null
